$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5 YEAR FORECAST")

# D6 had a hard-coded value (156000*103%) where it should reference B6 instead.
$ws.Range("D6").Formula = "=B6*103%"

$excel.Calculate()
